$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts position_id/tax_src_id/tax_dest_id right)
$ws.Columns("B:B").Insert()

# New header for the inserted column
$ws.Range("B1").Value = "_requirements"

# Populate the "_requirements" notes for the relevant rows
$ws.Range("B6").Value = "l10n_it_reverse_charge"
$ws.Range("B7").Value = "l10n_it_reverse_charge"
$ws.Range("B8").Value = "l10n_it_split_payment"
$ws.Range("B9").Value = "l10n_it_dichiarazione_intento or l10n_it_lettera_intento "

# Column widths (A stays, B is new/wide, C/D/E keep old B/C/D widths)
$ws.Columns("A:A").ColumnWidth = 19.77
$ws.Columns("B:B").ColumnWidth = 45.62
$ws.Columns("C:C").ColumnWidth = 17.83
$ws.Columns("D:D").ColumnWidth = 13.37
$ws.Columns("E:E").ColumnWidth = 17.27

# Selection, matching the final file
$ws.Range("B10").Select()

Write-Host "done"
